$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I
$ws.Range("I1").Value = "Other found locations"

# Row 2 - Kapil Goyal et al. (extra spacing in author list), elsevier/PMC tag
$ws.Range("E2").Value = "[Kapil%Goyal%kapilgoyalpgi@gmail.com%1,    Poonam%Chauhan%chauhan.poonamk@gmail.com%1,    Komal%Chhikara%dimpichhikara@gmail.com%1,    Parakriti%Gupta%parakritii@gmail.com%1,    Mini P.%Singh%minipsingh@gmail.com%1]"
$ws.Range("I2").Value = "_PMC_elsevier"

# Row 3 - Mamun / Griffiths
$ws.Range("E3").Value = "[Mohammed A.%Mamun%NULL%1,    Mark D.%Griffiths%NULL%6]"
$ws.Range("I3").Value = "_PMC_elsevier"

# Row 4 - Swapnajeet Sahoo et al.
$ws.Range("E4").Value = "[Swapnajeet%Sahoo%NULL%1,    Seema%Rani%NULL%1,    Shaheena%Parveen%NULL%1,    Ajay%Pal Singh%NULL%1,    Aseem%Mehra%NULL%1,    Subho%Chakrabarti%NULL%1,    Sandeep%Grover%NULL%1,    Cheering%Tandup%NULL%2,    Cheering%Tandup%NULL%0]"
$ws.Range("I4").Value = "_PMC_elsevier"

# Row 5 - Department of Health et al. (CORE record); ID/ID Format now "not found"/"N/A"
$ws.Range("E5").Value = "[Department%of Health%coreGivesNoEmail%1,   Department%of Health of the Government of the Hong Kong Special Administrative Region%coreGivesNoEmail%1,   US%Department of Health and Human Services`u{2014}Public Health Service%coreGivesNoEmail%1,   WHO%European Ministerial Conference on Mental Health%coreGivesNoEmail%1,   WHO%Website%coreGivesNoEmail%1]"
$ws.Range("F5").Value = "not found"
$ws.Range("G5").Value = "N/A"

# Row 6 - Chau / Cheung / Yip (CORE record); ID/ID Format now "not found"/"N/A"
$ws.Range("E6").Value = "[Chau%PH%coreGivesNoEmail%1,   Cheung%YT%coreGivesNoEmail%1,   Yip%PSF%coreGivesNoEmail%1]"
$ws.Range("F6").Value = "not found"
$ws.Range("G6").Value = "N/A"

# Row 7 - Chien-Cheng Huang et al.
$ws.Range("E7").Value = "[Chien-Cheng%Huang%NULL%1,    David Hung-Tsang%Yen%hjyen@vghtpe.gov.tw%1,    Hsien-Hao%Huang%NULL%1,    Wei-Fong%Kao%NULL%1,    Lee-Min%Wang%NULL%1,    Chun-I%Huang%NULL%1,    Chen-Hsen%Lee%NULL%1]"
$ws.Range("I7").Value = "_PMC_elsevier"

# Row 9 - Olaoluwa Okusaga et al.
$ws.Range("E9").Value = "[Olaoluwa%Okusaga%NULL%1,    Robert H.%Yolken%NULL%1,    Patricia%Langenberg%NULL%1,    Manana%Lapidus%NULL%1,    Timothy A.%Arling%NULL%1,    Faith B.%Dickerson%NULL%1,    Debra A.%Scrandis%NULL%1,    Emily%Severance%NULL%1,    Johanna A.%Cabassa%NULL%1,    Theodora%Balis%NULL%1,    Teodor T.%Postolache%NULL%1]"
$ws.Range("I9").Value = "_PMC_elsevier"

# Row 10 - Karine Kahil et al.
$ws.Range("E10").Value = "[Karine%Kahil%NULL%1,    Mohamad Ali%Cheaito%NULL%1,    Rawad%El Hayek%NULL%1,    Marwa%Nofal%NULL%1,    Sarah%El Halabi%NULL%1,    Kundadak Ganesh%Kudva%NULL%1,    Victor%Pereira-Sanchez%NULL%1,    Samer%El Hayek%NULL%1]"
$ws.Range("I10").Value = "_PMC_elsevier"

# Rows 5, 6, 8 also gained an (empty) "Other found locations" entry in the diff.
$ws.Range("I5").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("I8").Value = ""
